# Updated SO LUI TC upto Invoice creation
# Insert a new "ProductTypeIndex" column at the front of the AddLine sheet,
# numbering the existing product-type rows 1, 2, ...

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddLine")

# Make room for the new first column; existing columns A/B shift to B/C.
$ws.Columns("A").Insert()

$ws.Range("A1").Value = "ProductTypeIndex"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Match column sizing to content, same as the other bestFit columns in this sheet.
$ws.Columns("A").ColumnWidth = 16.6666666666667

$null = $ws.Range("A5").Select()
